$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BSD_T_length ")
$ws.Rows.Item(7).Insert()
